# Order row 4 (2022-09-17 19:35:33 / @dimchxn) has had its payment confirmed,
# so its status moves from "handling" to "payed".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J4").Value = "payed"
